# Update "想去人数" (F column) figures across the sheets that were
# refreshed by the gh-pages data generation run (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 113
$ws1.Range("F3").Value  = 183
$ws1.Range("F4").Value  = 424
$ws1.Range("F5").Value  = 200
$ws1.Range("F6").Value  = 135
$ws1.Range("F7").Value  = 1181
$ws1.Range("F10").Value = 54
$ws1.Range("F12").Value = 379
$ws1.Range("F13").Value = 407
$ws1.Range("F16").Value = 729
$ws1.Range("F19").Value = 1018
$ws1.Range("F20").Value = 474
$ws1.Range("F21").Value = 271
$ws1.Range("F22").Value = 88
$ws1.Range("F25").Value = 44
$ws1.Range("F26").Value = 477

# --- Sheet "演出" ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 150

# --- Sheet "全部类型" ------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 114
$ws4.Range("F5").Value  = 183
$ws4.Range("F6").Value  = 424
$ws4.Range("F7").Value  = 200
$ws4.Range("F9").Value  = 1181
$ws4.Range("F13").Value = 54
$ws4.Range("F17").Value = 379
$ws4.Range("F20").Value = 407
$ws4.Range("F23").Value = 729
$ws4.Range("F26").Value = 1018
$ws4.Range("F27").Value = 474
$ws4.Range("F30").Value = 271
$ws4.Range("F31").Value = 88
$ws4.Range("F34").Value = 150
$ws4.Range("F36").Value = 44
$ws4.Range("F38").Value = 477

$wb.Save()
